$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set attendance values from 0 to 1 for the specified cells
$ws.Range("G3").Value = 1
$ws.Range("H3").Value = 1

$ws.Range("D4").Value = 1
$ws.Range("E4").Value = 1

$ws.Range("D5").Value = 1
$ws.Range("E5").Value = 1

$ws.Range("G6").Value = 1
$ws.Range("H6").Value = 1

$ws.Range("H7").Value = 1
$ws.Range("H8").Value = 1
$ws.Range("H9").Value = 1
$ws.Range("H10").Value = 1
$ws.Range("H11").Value = 1
$ws.Range("H12").Value = 1
$ws.Range("H13").Value = 1
$ws.Range("H14").Value = 1
$ws.Range("H15").Value = 1
$ws.Range("H16").Value = 1

$ws.Range("D17").Value = 1
$ws.Range("E17").Value = 1

$ws.Range("H18").Value = 1
